$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.227.98'
$ws.Range('E2').Value = '  -3.40%  '
$ws.Range('D3').Value = '3.365.24'
$ws.Range('E3').Value = '  -4.18%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.39'
$ws.Range('E5').Value = '  -3.56%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '124.89'
$ws.Range('E6').Value = '  -7.25%  '
$ws.Range('D8').Value = '3.362.69'
$ws.Range('E8').Value = '  -4.24%  '
$ws.Range('E9').Value = '  -3.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.21'
$ws.Range('E10').Value = '  -5.62%  '
$ws.Range('E11').Value = '  -4.64%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.373'
$ws.Range('E12').Value = '  -3.96%  '
$ws.Range('D13').Value = '3.940.48'
$ws.Range('E13').Value = '  -4.12%  '
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D15').Value = '3.368.73'
$ws.Range('E15').Value = '  -4.07%  '
$ws.Range('E16').Value = '  -6.40%  '
$ws.Range('D17').Value = '62.281.78'
$ws.Range('E17').Value = '  -3.29%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '24.28'
$ws.Range('E18').Value = '  -5.87%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.18'
$ws.Range('E19').Value = '  -8.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.61'
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.03'
$ws.Range('E21').Value = '  -4.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '371.19'
$ws.Range('E22').Value = '  -6.38%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.551'
$ws.Range('E23').Value = '  -4.80%  '
$ws.Range('D24').Value = '3.500.88'
$ws.Range('E24').Value = '  -4.12%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '71.21'
$ws.Range('E26').Value = '  -4.65%  '
$ws.Range('E27').Value = '  -11.11%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.87'
$ws.Range('E29').Value = '  -7.52%  '
$ws.Range('E30').Value = '  -7.24%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.74'
$ws.Range('E31').Value = '  -6.24%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.38'
$ws.Range('E33').Value = '  -7.20%  '
$ws.Range('D34').Value = '3.394.48'
$ws.Range('E34').Value = '  -4.11%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.147'
$ws.Range('E35').Value = '  -6.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '22.56'
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.15'
$ws.Range('E37').Value = '  -4.00%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '165.63'
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.60'
$ws.Range('E39').Value = '  -5.33%  '
$ws.Range('E40').Value = '  -6.15%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0747'
$ws.Range('E41').Value = '  -5.29%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.761'
$ws.Range('E43').Value = '  -6.22%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '41.35'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.21'
$ws.Range('E45').Value = '  -5.51%  '
$ws.Range('E46').Value = '  -9.64%  '
$ws.Range('E47').Value = '  -8.24%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.06'
$ws.Range('E48').Value = '  -9.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.56'
$ws.Range('E49').Value = '  -3.65%  '
$ws.Range('D50').Value = '2.234.44'
$ws.Range('E50').Value = '  -6.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.832'
$ws.Range('E51').Value = '  -8.46%  '
